$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -0.043648266759025681
$ws.Range("B1").Value = 0.043648265827040789

$ws.Range("A2").Value = -0.023995201930941928
$ws.Range("B2").Value = 0.023995200960853242

$ws.Range("A3").Value = 0.051862731487112269
$ws.Range("B3").Value = -0.051862732454726605

$ws.Range("A4").Value = -0.0035167632268820709
$ws.Range("B4").Value = 0.0035167622195831374
